# Generate Report for Handoff
#
# "b.md" has progressed from "Handed back: in sync with en-US" to
# "Ready for handoff" in both locales. New handoff xliff files were
# generated (b.63290e5768f688058c7b37413b0a5c26c308f864.*.xlf), timestamps
# were refreshed, and the handback file version mismatch is now reported
# in the Error Detail column.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet (row 3 = b.md) ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-18 02:32:16"

# ---- zh-cn sheet (row 3 = b.md) ----
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-18 02:32:11"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1e77e170bf964edef9e186598b054b41226880e7/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/710d155568aff15f0d7c221de1f4e896328aa88c/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet (row 3 = b.md) ----
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-18 02:32:16"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1e77e170bf964edef9e186598b054b41226880e7/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/710d155568aff15f0d7c221de1f4e896328aa88c/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.17
